$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("BK1").Value = 0.8817262584677813
$ws.Range("BO1").Value = 0.82195294530689345
$ws.Range("BP1").Value = 0.99368594228649421
$ws.Range("A2").Value = 0.81141869715087767
$ws.Range("F2").Value = 0.76067503995300867
$ws.Range("D3").Value = 0.85335361535241927
$ws.Range("E3").Value = 0.62833478179923796
$ws.Range("AZ3").Value = 0.93360129861531183
$ws.Range("B4").Value = 0.84147641073279433
$ws.Range("AU4").Value = 0.56376327330255249
$ws.Range("D5").Value = 0.74796795818088646
$ws.Range("D6").Value = 0.82451270230407603
$ws.Range("E6").Value = 0.70042259237232329
$ws.Range("AU6").Value = 0.70075847768704813
$ws.Range("AB7").Value = 0.95935099857498263
$ws.Range("AZ7").Value = 0.81037237640020154
$ws.Range("I8").Value = 0.95214184895496745
$ws.Range("J8").Value = 0.85329629603597212
$ws.Range("E9").Value = 0.97325543371643553
$ws.Range("G9").Value = 0.92965445316844098
$ws.Range("I10").Value = 0.9247642569851946
$ws.Range("K10").Value = 0.82290413753519565
$ws.Range("AS10").Value = 0.97520670135663523
$ws.Range("I11").Value = 0.77817373257446221
$ws.Range("K12").Value = 0.93362013275440947
$ws.Range("N12").Value = 0.99623185360154487
$ws.Range("L13").Value = 0.73204105912723849
$ws.Range("X13").Value = 0.98677634027451511
$ws.Range("AU13").Value = 0.8736229842916704
$ws.Range("N15").Value = 0.76267432623226727
$ws.Range("P15").Value = 0.98408659895555028
$ws.Range("Q15").Value = 0.70556368150197435
$ws.Range("N16").Value = 0.71384192436185812
$ws.Range("Q16").Value = 0.95066544018998533
$ws.Range("R16").Value = 0.81477449996530749
$ws.Range("R17").Value = 0.92914495920153417
$ws.Range("S17").Value = 0.9588365168981926
$ws.Range("T18").Value = 0.86725325226334782
$ws.Range("R19").Value = 0.99662504051544565
$ws.Range("T19").Value = 0.72385992267721622
$ws.Range("S21").Value = 0.88177247810034798
$ws.Range("T21").Value = 0.80381642507008277
$ws.Range("T22").Value = 0.98070466437531179
$ws.Range("W22").Value = 0.81125845043955525
$ws.Range("X22").Value = 0.6983893963839376
$ws.Range("U23").Value = 0.97702296729724392
$ws.Range("Y23").Value = 0.98103261321475976
$ws.Range("BG23").Value = 0.95582651352106796
$ws.Range("F24").Value = 0.91499647160225384
$ws.Range("Z24").Value = 0.91533833413724452
$ws.Range("X25").Value = 0.82719837854134992
$ws.Range("Z25").Value = 0.74817646687877948
$ws.Range("AA25").Value = 0.95855904472438713
$ws.Range("AA26").Value = 0.91832222191465918
$ws.Range("AB27").Value = 0.73922659801981161
$ws.Range("AC27").Value = 0.90304683500501381
$ws.Range("Z28").Value = 0.63530877990200874
$ws.Range("AW28").Value = 0.890389370307068
$ws.Range("BN28").Value = 0.79717999755790059
$ws.Range("AS29").Value = 0.948377982036542
$ws.Range("BE29").Value = 0.85509856819773911
$ws.Range("AD31").Value = 0.83826750137477912
$ws.Range("AK31").Value = 0.90730852294855069
$ws.Range("BL31").Value = 0.76519878660715546
$ws.Range("AD32").Value = 0.84495181152185683
$ws.Range("AH32").Value = 0.84069691114238476
$ws.Range("AF33").Value = 0.9704683646753689
$ws.Range("AH33").Value = 0.98077147783880436
$ws.Range("BF34").Value = 0.9980435461816024
$ws.Range("AG35").Value = 0.96516680252661891
$ws.Range("AH35").Value = 0.77908455182382397
$ws.Range("AK35").Value = 0.80653649693079044
$ws.Range("AM35").Value = 0.69019539518642259
$ws.Range("BK35").Value = 0.59601740738770004
$ws.Range("AH37").Value = 0.85649197025872259
$ws.Range("AL37").Value = 0.80066069797922279
$ws.Range("AM37").Value = 0.76555961995660582
$ws.Range("AS37").Value = 0.9913343166452111
$ws.Range("BP37").Value = 0.9774447191557456
$ws.Range("O38").Value = 0.99735300844635044
$ws.Range("AJ38").Value = 0.85022825024899074
$ws.Range("BN38").Value = 0.9158831950033498
$ws.Range("N39").Value = 0.6729578013823605
$ws.Range("U39").Value = 0.99557017778433954
$ws.Range("AL39").Value = 0.82192479343089986
$ws.Range("C40").Value = 0.78033705456683056
$ws.Range("AQ41").Value = 0.75589723875193027
$ws.Range("AN42").Value = 0.77747661736437546
$ws.Range("AP43").Value = 0.87121105760178952
$ws.Range("AR43").Value = 0.76910272789316037
$ws.Range("AP44").Value = 0.60021220881692328
$ws.Range("AS44").Value = 0.9793543621826124
$ws.Range("AV44").Value = 0.82873315663735581
$ws.Range("AR46").Value = 0.81871157541574857
$ws.Range("AU46").Value = 0.84772381236824312
$ws.Range("AV46").Value = 0.79402532996455055
$ws.Range("P47").Value = 0.59888342618317658
$ws.Range("BH47").Value = 0.81924629269305638
$ws.Range("AQ48").Value = 0.87400447370153167
$ws.Range("AW48").Value = 0.84337624382863496
$ws.Range("AX49").Value = 0.84554834039346327
$ws.Range("Y50").Value = 0.74980632914757217
$ws.Range("AZ50").Value = 0.70876121068532216
$ws.Range("AZ51").Value = 0.80756672510486316
$ws.Range("BA51").Value = 0.65383335374914753
$ws.Range("BB52").Value = 0.92376814939686169
$ws.Range("BB53").Value = 0.90678409729045639
$ws.Range("BC54").Value = 0.67719085810407265
$ws.Range("BA55").Value = 0.67832980737008275
$ws.Range("BD55").Value = 0.98589722370779853
$ws.Range("BB56").Value = 0.81775285080038551
$ws.Range("BE56").Value = 0.90394598461962306
$ws.Range("BF56").Value = 0.75452138272592861
$ws.Range("BC57").Value = 0.98276417272664707
$ws.Range("BG57").Value = 0.86289279800021235
$ws.Range("BE58").Value = 0.8399389466000915
$ws.Range("BF60").Value = 0.90239309461569639
$ws.Range("BG60").Value = 0.92976189514736896
$ws.Range("K61").Value = 0.9926188906461807
$ws.Range("BG61").Value = 0.78287894589644047
$ws.Range("BJ61").Value = 0.95086203259462865
$ws.Range("BH62").Value = 0.74308912557713369
$ws.Range("BK62").Value = 0.5981268151773349
$ws.Range("BL62").Value = 0.87026542875572532
$ws.Range("AO63").Value = 0.838160941438719
$ws.Range("BI63").Value = 0.88260623148729023
$ws.Range("E64").Value = 0.85285772075346999
$ws.Range("AJ64").Value = 0.99676411110317342
$ws.Range("AS64").Value = 0.78463753762699784
$ws.Range("BK64").Value = 0.95711148957232406
$ws.Range("BN64").Value = 0.82361379638281473
$ws.Range("BN65").Value = 0.90882405163919755
$ws.Range("BO65").Value = 0.62547958478636223
$ws.Range("C66").Value = 0.74708452325388208
$ws.Range("AI67").Value = 0.85641141449778502
$ws.Range("BP67").Value = 0.97533605176405991
